# edit.ps1 - applies the OOXML diff:
#   - slide 3: resize "Rectangle 5", reposition/resize "Straight Arrow Connector 38"
#   - slide 7: update the title text
#   - slide 8: reposition title placeholder and picture
#   - slide 9: reposition title placeholder and picture
#
# Note: PowerPoint COM exposes Left/Top/Width/Height in points (1 pt = 12700 EMU),
# and the host here stores them as single-precision floats, truncating (floor) when
# converting back to EMU. The literals below were chosen so that
# floor(float32(value) * 12700) reproduces the exact target EMU from the diff.

$p = $ppt.ActivePresentation

# --- Slide 3 -------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$rect5 = $null
$conn38 = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $sh = $s3.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 5") { $rect5 = $sh }
    if ($sh.Name -eq "Straight Arrow Connector 38") { $conn38 = $sh }
}

# ext cy: 432729 -> 489496 (off/ cx unchanged)
$rect5.Height = 38.54299545288086

# off y: 3936782 -> 4003653 (off x unchanged)
# ext cx: 10070 -> 13468
# ext cy: 650259 -> 587295
$conn38.Top = 315.248291015625
$conn38.Width = 1.0604724884033203
$conn38.Height = 46.24370193481445

# --- Slide 7 ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$title7 = $s7.Shapes.Item(1)
$title7.TextFrame.TextRange.Text = "Plots of Count of Total Investment and Average Investment "

# --- Slide 8 ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)

$title8 = $null
$pic8 = $null
for ($i = 1; $i -le $s8.Shapes.Count; $i++) {
    $sh = $s8.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") { $title8 = $sh }
    if ($sh.Name -eq "Picture 1") { $pic8 = $sh }
}

# off x: 1136468 -> 1232262, off y: 70529 -> 627877 (ext unchanged)
$title8.Left = 97.02851104736328
$title8.Top = 49.43913650512695

# off x: 1514611 -> 1436233, off y: 1113392 -> 1592363 (ext unchanged)
$pic8.Left = 113.08921813964844
$pic8.Top = 125.38291931152344

# --- Slide 9 ---------------------------------------------------------------
$s9 = $p.Slides.Item(9)

$title9 = $null
$pic9 = $null
for ($i = 1; $i -le $s9.Shapes.Count; $i++) {
    $sh = $s9.Shapes.Item($i)
    if ($sh.Name -eq "Title 1") { $title9 = $sh }
    if ($sh.Name -eq "Picture 3") { $pic9 = $sh }
}

# off x: 1182651 -> 1226194, off y: 230744 -> 822927 (ext unchanged)
$title9.Left = 96.55071258544922
$title9.Top = 64.79740905761719

# off x unchanged (1841093), off y: 1323702 -> 1820415 (ext unchanged)
$pic9.Top = 143.3397674560547
